$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("test results")

# Fix F13: was stored as text "7.67", convert it to a real number.
$ws.Cells.Item(13, 6).Value = 7.67

# Add new row 14 with the test result data.
$ws.Cells.Item(14, 1).Value = 13
$ws.Cells.Item(14, 2).Value = "asdfasf"
$ws.Cells.Item(14, 3).Value = "19-01-202500:00:53"
$ws.Cells.Item(14, 4).Value = "EN->PL"
$ws.Cells.Item(14, 5).Value = 240

# F14 and H14 are numeric/percent-looking strings that must stay as text,
# so force text format before assigning, then drop the format again so the
# cell ends up with no explicit style (matching rows 2-13).
$f14 = $ws.Cells.Item(14, 6)
$f14.NumberFormat = "@"
$f14.Value = "10.74"
$f14.ClearFormats()

$ws.Cells.Item(14, 7).Value = "(2/None)"

$h14 = $ws.Cells.Item(14, 8)
$h14.NumberFormat = "@"
$h14.Value = "50.00%"
$h14.ClearFormats()

$ws.Cells.Item(14, 9).Value = "Games Remastering"
